$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C-G (rows 2-17), per "Last Update 15-03-2024" diff.
# Values must remain stored as literal TEXT (matching the original file's
# convention of text-typed odds/percentages) rather than being auto-converted
# to numbers/percrentages by Excel, and the default (no explicit) cell style
# must be preserved.
$data = @{
    2 = @{ C = "1.9"; D = "5.1"; E = "79%"; F = "32%"; G = "2.32" }
    3 = @{ C = "1.4"; D = "6.4"; E = "90%"; F = "52%"; G = "2.74" }
    4 = @{ C = "1.9"; D = "5.6"; E = "79%"; F = "52%"; G = "2.84" }
    5 = @{ C = "2.3"; D = "5.6"; E = "79%"; F = "53%"; G = "2.89" }
    6 = @{ C = "1.8"; D = "5.4"; E = "79%"; F = "52%"; G = "2.68" }
    7 = @{ C = "2.0"; D = "5.1"; E = "100%"; F = "74%"; G = "3.42" }
    8 = @{ C = "1.7"; D = "4.2"; E = "79%"; F = "53%"; G = "3.00" }
    9 = @{ C = "1.6"; D = "4.5"; E = "58%"; F = "37%"; G = "2.11" }
    10 = @{ C = "2.2"; D = "2.5"; E = "52%"; F = "16%"; G = "1.74" }
    11 = @{ C = "1.9"; D = "5.5"; E = "79%"; F = "58%"; G = "2.89" }
    12 = @{ C = "1.9"; D = "5.1"; E = "68%"; F = "31%"; G = "1.89" }
    13 = @{ C = "2.6"; D = "4.5"; E = "68%"; F = "47%"; G = "2.21" }
    14 = @{ C = "1.7"; D = "5.5"; E = "63%"; F = "48%"; G = "2.47" }
    15 = @{ C = "2.3"; D = "5.4"; E = "95%"; F = "37%"; G = "2.58" }
    16 = @{ C = "1.6"; D = "4.8"; E = "53%"; F = "32%"; G = "1.89" }
    17 = @{ C = "2.1"; D = "4.1"; E = "79%"; F = "42%"; G = "2.53" }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in @("C", "D", "E", "F", "G")) {
        $cell = $ws.Range("$col$r")
        $cell.NumberFormat = "@"
        $cell.Value = $row[$col]
        $cell.Style = "Normal"
    }
}

Write-Host "Updated odds/percentage columns C:G for rows 2-17"
